# PressureSensorCalibration.xlsx - recalibrated sensor readings (femur broke, needs redesign)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated raw voltage readings (column A) from the new calibration run.
$ws.Range("A2").Value = 0.8407
$ws.Range("A3").Value = 2.1065
$ws.Range("A4").Value = 2.7713
$ws.Range("A5").Value = 3.4262
$ws.Range("A6").Value = 4.0518
$ws.Range("A7").Value = 4.6872
$ws.Range("A8").Value = 4.9022

# Min/max helper points (columns M/N) used for the trend overlay chart mirror A2/A8.
$ws.Range("M2").Value = 0.8407
$ws.Range("M3").Value = 4.9022

# New sample row showing the percent drift between the old and new voltage reading.
$ws.Range("F30").Formula = "=(2.1065-2.0821)/2.1065"

# Update the active cell/selection to reflect where editing left off.
[void]$ws.Range("M26").Select()
